# V1.0.3 - Add "checkMagCondition" field documentation row to the
# SwitchableRangedWeapon reference sheet (Misc category), matching the
# upstream commit that documents the new magazine-condition check property.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Write the new row's content (column order matters so that the new
# shared strings are appended to sharedStrings.xml in the same order as
# the reference edit: 杂项, checkMagCondition, 设置是否在开火时检查弹匣耐久) ---
$ws.Range("A11").Value = "杂项"
$ws.Range("B11").Value = "checkMagCondition"
$ws.Range("C11").Value = "为一个布尔值"
$ws.Range("D11").Value = "设置是否在开火时检查弹匣耐久"
$ws.Range("E11").Value = "/"
$ws.Range("F11").Value = "/"
$ws.Range("G11").Value = $true

# --- Formatting: centered text, boxed row (thin borders), matching the
# look of the existing table rows above it ---
$newRow = $ws.Range("A11:G11")
$newRow.HorizontalAlignment = -4108   # xlCenter
$newRow.VerticalAlignment = -4108     # xlCenter

# Give the whole new row a top and bottom edge (applied range-wide so all
# seven cells pick up one shared style).
$newRow.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$newRow.Borders.Item(9).LineStyle = 1   # xlEdgeBottom

# A11 and C11 additionally get a left edge (matching the column separator
# used for the "field name" / "category" columns elsewhere in the sheet),
# and G11 gets a right edge closing off the row.
$ws.Range("A11").Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$ws.Range("C11").Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$ws.Range("G11").Borders.Item(10).LineStyle = 1  # xlEdgeRight

# --- Update the active selection to the newly added row, like the author
# left it selected after editing ---
$null = $ws.Range("A11:G11").Select()
